$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B4 with a formula that evaluates to the commit hash text "514716e1" ---
$ws.Range("B4").Formula = '="514716e1"'
$ws.Range("B4").NumberFormat = "0.00E+00"
$ws.Range("B4").HorizontalAlignment = -4108

# --- Add new row 14: "No" / "Was version properly marked in NUGET configuration before commit?" ---
$ws.Range("A14").Value = "No"
$ws.Range("B14").Value = "Was version properly marked in NUGET configuration before commit?"
$ws.Range("B14:J14").HorizontalAlignment = -4108
$ws.Range("B14:J14").Merge()

# --- Add new row 15: "Yes" / "Has version in project configuration and NUGET config been updated?" ---
$ws.Range("A15").Value = "Yes"
$ws.Range("B15").Value = "Has version in project configuration and NUGET config been updated?"
$ws.Range("B15:J15").HorizontalAlignment = -4108
$ws.Range("B15:J15").Merge()

# --- Update selection to reflect the newly added row ---
$ws.Range("B14:J14").Select()
